$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add column O (Terms, copied from column A) and column P (numeric values) for rows 2-23 ---
$pValues = @{
    2  = 0.23076792961122999
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0.97300881941683603
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 0
}

for ($r = 2; $r -le 23; $r++) {
    $term = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 15).Value2 = $term
    $ws.Cells.Item($r, 16).Value2 = $pValues[$r]
}

# --- Register the generic "0.00" number format (as used on column R) via a scratch cell ---
# that is removed again afterwards without disturbing any real data (shift-left keeps all rows intact).
$ws.Range("Z1").NumberFormat = "0.00"
$ws.Range("Z1").Delete(-4159) | Out-Null

# --- R18 gets a more precise custom number format; the cell itself stays empty ---
$ws.Range("R18").NumberFormat = "0.0000"
$ws.Range("R18").Value = $null

# --- Update the selection / active cell to R18:S18 like in the saved file ---
$ws.Activate()
$ws.Range("R18:S18").Select()
